$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5605295187
$ws.Range("C2").Value = -224.73613661
$ws.Range("D2").Value = -225.29666612
$ws.Range("E2").Value = -224.5113559466

$ws.Range("B3").Value = -0.569162016
$ws.Range("C3").Value = -224.66421932
$ws.Range("D3").Value = -225.23338133
$ws.Range("E3").Value = -224.5113559466

$ws.Range("B4").Value = -0.5719302858
$ws.Range("C4").Value = -224.6448457
$ws.Range("D4").Value = -225.21677599
$ws.Range("E4").Value = -224.5113559466
